# Run 181 refresh of the optimisation result: the "Schedule" sheet gets the
# latest ON/OFF block summary, and the "Detailed" sheet drops its oldest
# (now stale) historical reading, shifting every remaining row up by one and
# refreshing the recomputed Price / Pump_Status values for the new run.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Schedule"
$ws2 = $wb.Worksheets.Item(2)   # "Detailed"

# --- Sheet1 "Schedule": update row 2, add new row 3 ---
$ws1.Cells.Item(2,1).Value2 = 46045.16666666666
$ws1.Cells.Item(2,3).Value2 = 12
$ws1.Cells.Item(2,4).Value2 = 45.36
$ws1.Cells.Item(2,5).Value2 = 812.9197342500003
$ws1.Cells.Item(2,6).Value2 = 17.92151089616403

$ws1.Cells.Item(3,1).Value2 = 46045.83333333334
$ws1.Cells.Item(3,2).Value2 = 46046
$ws1.Cells.Item(3,3).Value2 = 4
$ws1.Cells.Item(3,4).Value2 = 15.12
$ws1.Cells.Item(3,5).Value2 = 447.2168902500001
$ws1.Cells.Item(3,6).Value2 = 29.57783665674604

# Match date-style formatting of row 2 (A/B columns) on the newly added row 3
$ws1.Range("A3:B3").NumberFormat = $ws1.Range("A2:B2").NumberFormat

# --- Sheet2 "Detailed": drop the oldest historical reading, shift everything up ---
$ws2.Rows(2).Delete()

# The delete above shifts old row (r+1) data into row r for every remaining
# row, which already gets DateTime (A), Type (C) and Date (D) correct for
# most rows. Refresh the recomputed Price (B) / Pump_Status (E) values (and
# the couple of Type cells that don't line up after the shift) to match the
# new run's results.
$ws2.Cells.Item(2,5).Value2 = "OFF"
$ws2.Cells.Item(3,3).Value2 = "historical"
$ws2.Cells.Item(3,5).Value2 = "OFF"
$ws2.Cells.Item(4,5).Value2 = "OFF"
$ws2.Cells.Item(5,5).Value2 = "OFF"
$ws2.Cells.Item(6,2).Value2 = 57.06
$ws2.Cells.Item(6,5).Value2 = "OFF"
$ws2.Cells.Item(7,5).Value2 = "OFF"
$ws2.Cells.Item(8,2).Value2 = 57.06
$ws2.Cells.Item(8,5).Value2 = "OFF"
$ws2.Cells.Item(12,2).Value2 = 64.8901
$ws2.Cells.Item(13,2).Value2 = 64.8901
$ws2.Cells.Item(15,2).Value2 = 35.89162
$ws2.Cells.Item(16,2).Value2 = 9.992749999999999
$ws2.Cells.Item(17,2).Value2 = 4.45658
$ws2.Cells.Item(18,2).Value2 = 34.16907
$ws2.Cells.Item(20,2).Value2 = 3.53702
$ws2.Cells.Item(21,2).Value2 = 22.37639
$ws2.Cells.Item(24,2).Value2 = 35.70911
$ws2.Cells.Item(25,2).Value2 = 36.06
$ws2.Cells.Item(26,2).Value2 = 36.06
$ws2.Cells.Item(27,2).Value2 = 36.06
$ws2.Cells.Item(28,2).Value2 = 36.06
$ws2.Cells.Item(32,2).Value2 = 4.53109
$ws2.Cells.Item(33,2).Value2 = 5.42477
$ws2.Cells.Item(34,2).Value2 = -1.97213
$ws2.Cells.Item(35,2).Value2 = -4.50626
$ws2.Cells.Item(36,2).Value2 = 4.7946
$ws2.Cells.Item(37,2).Value2 = 44.28147
$ws2.Cells.Item(38,2).Value2 = 57.01493
$ws2.Cells.Item(39,2).Value2 = 57.04922
$ws2.Cells.Item(40,2).Value2 = 59.30893
$ws2.Cells.Item(41,2).Value2 = 59.38449
$ws2.Cells.Item(41,5).Value2 = "ON"
$ws2.Cells.Item(42,5).Value2 = "ON"
$ws2.Cells.Item(43,2).Value2 = 57.06
$ws2.Cells.Item(43,5).Value2 = "ON"
$ws2.Cells.Item(44,5).Value2 = "ON"
$ws2.Cells.Item(45,2).Value2 = 56.9895
$ws2.Cells.Item(45,5).Value2 = "ON"
$ws2.Cells.Item(46,5).Value2 = "ON"
$ws2.Cells.Item(47,5).Value2 = "ON"
$ws2.Cells.Item(48,5).Value2 = "ON"
